# Hardware_and_Software.docx edit:
#   "Updated with Kibana link" / "Updated with Kibana Link"
#
# Content-level change: in the "External API Hoiio used for sending SMS..."
# paragraph, the sentence was edited so the trailing "... The API link" and
# " is as follows" runs (which had the hidden "_GoBack" last-edit bookmark
# sitting between them) become one contiguous run/sentence, and Word's
# "_GoBack" bookmark (which always marks the position of the most recent
# edit) ends up anchored at the start of the "(3) Readme file ..." paragraph
# instead - i.e. where the author's cursor/next edit was.

$d = $word.ActiveDocument

# 1) Merge " ... The API link" + " is as follows" into a single run by
#    re-typing the (already textually contiguous) text in place, starting
#    right after "Hoiio" (which keeps its own run + surrounding proofErr
#    spell-check tags untouched). This collapses the trailing two runs
#    (and the "_GoBack" bookmark that sat between them) into one run,
#    matching the diff's run merge.
$rng = $d.Content
$found = $rng.Find.Execute(
    " used for sending SMS. The API link is as follows",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " used for sending SMS. The API link is as follows", 2)

if (-not $found) {
    throw "Could not find the ' used for sending SMS. The API link is as follows' text to update."
}

# 2) Re-anchor the hidden "_GoBack" bookmark at the start of the
#    "(3) Readme file ..." paragraph (that's where it ends up after the
#    edit above, per the authored diff).
$target = $d.Content
$found2 = $target.Find.Execute("(3) Readme file that indicates", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the '(3) Readme file that indicates ...' paragraph."
}
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target)
